# Re-conversion using pandoc after adjusting image spacing under fig 1
#
# Splits the two "caption + image" paragraphs (Figure 3 / Example Heat Map)
# into: a caption-text paragraph, a CaptionedFigure paragraph holding just
# the picture, and a new ImageCaption paragraph repeating the caption text
# under the image. Mirrors the standard pandoc docx output for a Markdown
# image-with-title, where the title text becomes an ImageCaption paragraph
# right after the figure.

$d = $word.ActiveDocument

function Get-PackageXml([string]$bodyInnerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
        'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
        'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
        'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
        'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
        '<w:body>' + $bodyInnerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphByText([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text.StartsWith($needle)) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Location 1: "Figure 3: R-loop formation correlates with ectopic
# deposition of H3K4me3" caption + inline picture.
# ---------------------------------------------------------------------
$p1 = Find-ParagraphByText("Figure 3: R-loop formation correlates with ectopic deposition of H3K4me3")

$drawing1 = '<w:r><w:drawing><wp:inline><wp:extent cx="1188720" cy="1419606" /><wp:effectExtent b="0" l="0" r="0" t="0" /><wp:docPr descr="Figure 3: R-loop formation correlates with ectopic deposition of H3K4me3" title="" id="27" name="Picture" /><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="https://oup.silverchair-cdn.com/oup/backfile/Content_public/Journal/nar/51/1/10.1093_nar_gkac1155/2/m_gkac1155fig3.jpeg?Expires=1678997233&amp;Signature=q8YRcMZ3FPyio7jb06-5rSE3nRYj0DYF3ARHtCdA5AkRkfgt0O7VR7C~9suLERQqxDWD8tbsaN0AdLAoiCoBHSWtlZ1D2u-JBlyG712Br9Uc~lTN26ZCww-c8UBNWoYGtySP31RjdPAjrkTKEyCndlA9Qkvqob2FZ7JKMd8DXbyVPlr3UPeSAfqvRQ1hGBzmcvxmut4oP5eGCJCA7M-A1OKGQdJn-p9PUyaVc2MzO12A7bxvSeajJNR~kPS3PqFSMJ-715YBHw2fe67ZOnCZ1Yu3iKcnKjy3NEQ8phFO-p1dot~0Z9083SNL6NgfAyBBNlmwDp-tTfULg1CPgEDMGw__&amp;Key-Pair-Id=APKAIE5G5CRDK6RD3PGA" id="28" name="Picture" /><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1" /></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId26" /><a:stretch><a:fillRect /></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0" /><a:ext cx="1188720" cy="1419606" /></a:xfrm><a:prstGeom prst="rect"><a:avLst /></a:prstGeom><a:noFill /><a:ln w="9525"><a:noFill /><a:headEnd /><a:tailEnd /></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>'

$body1 = '<w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr>' +
    '<w:r><w:rPr><w:bCs /><w:b /></w:rPr><w:t xml:space="preserve">Figure 3: R-loop formation correlates with ectopic deposition of</w:t></w:r>' +
    '<w:r><w:rPr><w:bCs /><w:b /></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:bCs /><w:b /></w:rPr><w:t xml:space="preserve">H3K4me3</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="CaptionedFigure" /></w:pPr>' + $drawing1 + '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ImageCaption" /></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Figure 3: R-loop formation correlates with ectopic deposition of</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">H3K4me3</w:t></w:r>' +
    '</w:p>'

$rng1 = $d.Range($p1.Range.Start, $p1.Range.End)
$rng1.InsertXML((Get-PackageXml($body1)))

# ---------------------------------------------------------------------
# Location 2: "Example Heat Map from StackOverflow" caption + inline
# picture.
# ---------------------------------------------------------------------
$p2 = Find-ParagraphByText("Example Heat Map from")

$drawing2 = '<w:r><w:drawing><wp:inline><wp:extent cx="5334000" cy="4574689" /><wp:effectExtent b="0" l="0" r="0" t="0" /><wp:docPr descr="Example Heat Map" title="" id="31" name="Picture" /><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="https://i.stack.imgur.com/EpF3I.png" id="32" name="Picture" /><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1" /></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId30" /><a:stretch><a:fillRect /></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0" /><a:ext cx="5334000" cy="4574689" /></a:xfrm><a:prstGeom prst="rect"><a:avLst /></a:prstGeom><a:noFill /><a:ln w="9525"><a:noFill /><a:headEnd /><a:tailEnd /></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>'

$body2 = '<w:p><w:pPr><w:pStyle w:val="BodyText" /></w:pPr>' +
    '<w:r><w:rPr><w:bCs /><w:b /></w:rPr><w:t xml:space="preserve">Example Heat Map from</w:t></w:r>' +
    '<w:r><w:rPr><w:bCs /><w:b /></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:hyperlink r:id="rId29"><w:r><w:rPr><w:rStyle w:val="Hyperlink" /><w:bCs /><w:b /></w:rPr><w:t xml:space="preserve">StackOverflow</w:t></w:r></w:hyperlink>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="CaptionedFigure" /></w:pPr>' + $drawing2 + '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ImageCaption" /></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Example Heat Map</w:t></w:r>' +
    '</w:p>'

$rng2 = $d.Range($p2.Range.Start, $p2.Range.End)
$rng2.InsertXML((Get-PackageXml($body2)))

Write-Output "Done. Paragraphs now: $($d.Paragraphs.Count)"
